$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.74005651473999
$ws.Range("B1").Value = 2.359943628311157
$ws.Range("C1").Value = 2.561846256256104
$ws.Range("D1").Value = 3.302919387817383
$ws.Range("E1").Value = 1.481994032859802
